# Applies the "way to test best_strategies" commit:
#  - rename Sheet1 -> visual, Sheet3 -> best strategy compiler
#  - move the workbook window
#  - scroll Sheet2 down a bit (topLeftCell = A3)
#  - replace Sheet3's contents with a table of (n, count, l, r) params
#    and a CONCATENATE formula building a "go playGames(...)" string

$wb = $excel.ActiveWorkbook

$wsVisual = $wb.Worksheets.Item(1)
$wsVisual.Name = "visual"

$wsCompiler = $wb.Worksheets.Item(3)
$wsCompiler.Name = "best strategy compiler"

$wsScore = $wb.Worksheets.Item(2)
$wsScore.Application.ActiveWindow.ScrollRow = 3

# --- clear out the old Sheet3 data (C7:C94) ---
$wsCompiler.Cells.Clear()

# --- new data table, rows 3..52, columns A..D ---
$data = @(
    @(7,1,2,4),
    @(6,1,5,3),
    @(7,3,4,6),
    @(7,1,3,4),
    @(4,1,2,5),
    @(6,2,4,5),
    @(6,1,3,5),
    @(7,2,5,8),
    @(7,1,4,7),
    @(7,2,4,6),
    @(7,3,6,7),
    @(5,1,2,6),
    @(6,2,3,6),
    @(7,2,3,3),
    @(7,1,6,4),
    @(7,2,5,5),
    @(6,2,4,3),
    @(6,2,2,5),
    @(7,2,5,3),
    @(5,2,5,6),
    @(7,1,5,5),
    @(7,2,3,5),
    @(6,1,2,7),
    @(5,1,3,4),
    @(6,2,4,7),
    @(7,3,3,8),
    @(5,1,2,7),
    @(7,2,2,7),
    @(5,2,3,7),
    @(5,2,3,6),
    @(6,1,3,7),
    @(5,1,5,4),
    @(6,1,5,8),
    @(7,2,4,7),
    @(5,2,3,4),
    @(5,2,3,5),
    @(6,1,4,4),
    @(7,2,6,8),
    @(7,1,5,7),
    @(5,1,3,8),
    @(7,1,3,7),
    @(7,2,4,8),
    @(5,1,4,5),
    @(7,1,2,5),
    @(7,3,3,6),
    @(5,1,2,4),
    @(5,1,4,4),
    @(5,1,3,5),
    @(6,1,5,6),
    @(7,3,2,5)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $wsCompiler.Cells.Item($row, 1).Value = $vals[0]
    $wsCompiler.Cells.Item($row, 2).Value = $vals[1]
    $wsCompiler.Cells.Item($row, 3).Value = $vals[2]
    $wsCompiler.Cells.Item($row, 4).Value = $vals[3]
    $wsCompiler.Cells.Item($row, 5).Formula = "=CONCATENATE(""go playGames("",A$row,"","",B$row,"","",C$row,"","",D$row,"", 10, TRUE)"")"
}

$wsCompiler.Columns.Item(5).ColumnWidth = 42.83203125

$wsCompiler.Range("E3:E52").Select()

$excel.ActiveWindow.Left = 20460
$excel.ActiveWindow.Top = 4300
